# Automatic update of files.
# This script swaps the data of row pairs (8,9), (11,12) and (13,14)
# on the active sheet, reproducing the re-ordering seen in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-CellValues($Worksheet, $RowA, $RowB, $Columns) {
    foreach ($col in $Columns) {
        $rangeA = $Worksheet.Range("$col$RowA")
        $rangeB = $Worksheet.Range("$col$RowB")

        $valA = $rangeA.Value()
        $valB = $rangeB.Value()

        $rangeA.Value = $valB
        $rangeB.Value = $valA
    }
}

$swapColumns = @("A", "B", "E", "F", "G", "H", "Q", "R", "Z", "AB")
$swapColumnsPartial = @("A", "Q", "R", "Z", "AB")

# Rows 8 and 9 fully swap their species/record data.
Swap-CellValues $ws 8 9 $swapColumns

# Rows 11 and 12 fully swap their species/record data.
Swap-CellValues $ws 11 12 $swapColumns

# Rows 13 and 14 only swap id/coordinate/time data (species data is identical).
Swap-CellValues $ws 13 14 $swapColumnsPartial
